# Apply the "Holden scheme" update described by the commit.
#
# Net effect on the worksheet:
#  - Rows 16-19 (HKL 14-17), previously labelled with the
#    "HexGrid-90degTilt*degRes" series, are relabelled "Holden2.5/5/10/15".
#  - Four new rows (20-23, HKL 18-21) are appended carrying the
#    "HexGrid-90degTilt2.5/5/10/15degRes" labels that used to live on rows 16-19,
#    each filled with 1's across C:T exactly like every other data row.
#  - The bracketed-plane header row (row 2, C2:J2) is reordered.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Relabel existing rows 16-19 -> Holden series ---------------------------
$ws.Range("B16").Value = "Holden2.5"
$ws.Range("B17").Value = "Holden5"
$ws.Range("B18").Value = "Holden10"
$ws.Range("B19").Value = "Holden15"

# --- Reorder the [h,k,l]-plane headers on row 2 (C2:J2) ---------------------
$ws.Range("C2").Value = "[3, 1, 0]"
$ws.Range("D2").Value = "[2, 2, 2]"
$ws.Range("E2").Value = "[1, 1, 0]"
$ws.Range("F2").Value = "[3, 2, 1]"
$ws.Range("G2").Value = "[4, 0, 0]"
$ws.Range("H2").Value = "[2, 1, 1]"
$ws.Range("I2").Value = "[2, 0, 0]"
$ws.Range("J2").Value = "[2, 2, 0]"

# --- Append the four new Holden data rows (20-23) ---------------------------
# Seed the new rows by copying the formatting (bold/border/centered column-A
# style, plain data-cell style elsewhere) from rows 16-19, then overwrite the
# values with the correct HKL index / label / data.
$ws.Range("A16:T19").Copy($ws.Range("A20:T23"))

$dataCols = @("C","D","E","F","G","H","I","J","K","L","M","N","O","P","Q","R","S","T")

$newRows = @(
    @{ Row = 20; Hkl = 18; Label = "HexGrid-90degTilt2.5degRes" },
    @{ Row = 21; Hkl = 19; Label = "HexGrid-90degTilt5degRes" },
    @{ Row = 22; Hkl = 20; Label = "HexGrid-90degTilt10degRes" },
    @{ Row = 23; Hkl = 21; Label = "HexGrid-90degTilt15degRes" }
)

foreach ($nr in $newRows) {
    $r = $nr.Row
    $ws.Range("A$r").Value = $nr.Hkl
    $ws.Range("B$r").Value = $nr.Label
    foreach ($c in $dataCols) {
        $ws.Range("$c$r").Value = 1
    }
}

# --- Tidy up the used range: drop the now-unused U:AD columns ---------------
$ws.Range("U1:AD2").Clear()
